# Updates the crypto price/volume table (columns D and E) with refreshed
# quote data, matching the upstream GitHub Actions scrape commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.922.59'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.630.22'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.90'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.45%  '
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0608'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0881'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').Value = '1.862.08'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '1.652.81'
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.85'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').Value = '27.928.62'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '228.03'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.64'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.35'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  -3.00%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '154.64'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.38'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').Value = '1.419.76'
$ws.Range('E33').Value = '  +1.37%  '
$ws.Range('E34').Value = '  +1.11%  '
$ws.Range('E35').Value = '  +2.93%  '
$ws.Range('E36').Value = '  -1.28%  '
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('E38').Value = '  -0.77%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.555'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('E40').Value = '  -1.42%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.30%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '65.79'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('E43').Value = '  -1.33%  '
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('D45').Value = '1.771.01'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('E46').Value = '  -3.78%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '88.67'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.57'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('E51').Value = '  +0.04%  '
